$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.156.06'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').Value = '1.614.19'
$ws.Range('E4').Value = '  -0.44%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '213.26'
$cell.ClearFormats()
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('E6').Value = '  -0.44%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.482'
$cell.ClearFormats()
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  +1.75%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '18.44'
$cell.ClearFormats()
$ws.Range('E10').Value = '  +3.52%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '1.839.15'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').Value = '1.612.57'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').Value = '26.163.88'
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '60.82'
$cell.ClearFormats()
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('E19').Value = '  -0.43%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '198.61'
$cell.ClearFormats()
$ws.Range('E20').Value = '  +4.98%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '4.27'
$cell.ClearFormats()
$ws.Range('E21').Value = '  +2.63%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '9.50'
$cell.ClearFormats()
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('E23').Value = '  +1.71%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '142.62'
$cell.ClearFormats()
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  -0.43%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '15.22'
$cell.ClearFormats()
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('E29').Value = '  +0.40%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.17'
$cell.ClearFormats()
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('E31').Value = '  +3.46%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.15'
$cell.ClearFormats()
$ws.Range('E32').Value = '  +2.57%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.04'
$cell.ClearFormats()
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('D36').Value = '1.108.18'
$ws.Range('E36').Value = '  +0.61%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '0.0153'
$cell.ClearFormats()
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('E38').Value = '  -0.51%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.507'
$cell.ClearFormats()
$ws.Range('E39').Value = '  +2.65%  '
$ws.Range('E40').Value = '  -0.92%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.791'
$cell.ClearFormats()
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('E42').Value = '  +7.86%  '
$ws.Range('D43').Value = '1.750.49'
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('E44').Value = '  +1.22%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '93.20'
$cell.ClearFormats()
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('D46').Value = '0.0₆0109'
$ws.Range('E46').Value = '  +7.90%  '
$ws.Range('E47').Value = '  +8.24%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '54.02'
$cell.ClearFormats()
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('E49').Value = '  +0.17%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.410'
$cell.ClearFormats()
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  -0.28%  '
